$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @{Row=1006; B='lab.mod.label'; C='Mody'; Height=0}
    @{Row=1007; B='lab.mod.title'; C='Mody'; Height=0}
    @{Row=1008; B='lab.mod.filter.title'; C='Filtr modů'; Height=0}
    @{Row=1009; B='lab.mod.button.create'; C='Nový mod'; Height=0}
    @{Row=1010; B='lab.mod.table.name'; C='Název'; Height=0}
    @{Row=1011; B='lab.mod.table.vendor'; C='Výrobce'; Height=0}
    @{Row=1012; B='lab.mod.table.footer.label'; C='Počet modů [{{data.total}}]'; Height=0}
    @{Row=1013; B='lab.mod.context.menu'; C='Mod [{{data.name}}]'; Height=0}
    @{Row=1014; B='lab.mod.preview'; C='Náhled modu'; Height=0}
    @{Row=1015; B='lab.mod.button.edit'; C='Upravit mod'; Height=0}
    @{Row=1016; B='lab.mod.button.delete'; C='Odstranit mod'; Height=0}
    @{Row=1017; B='lab.mod.button.delete.confirm.title'; C='Odstranit mod'; Height=0}
    @{Row=1018; B='lab.mod.button.delete.confirm'; C='Opravdu si přejete odstranit vybraný mod? Bude smazána velké množství dat, která jsou spojená s jeho využitím ve vapování (případně jinde). Použijte s rozvahou, poněvadž není cesty zpět.'; Height=30}
    @{Row=1019; B='lab.mod.button.delete.confirm.ok'; C='Odstranit mod'; Height=0}
    @{Row=1020; B='lab.mod.deleted.success'; C='Mod [{{data.name}}] byl úspěšně odstraněn.'; Height=0}
    @{Row=1021; B='lab.mod.index.label'; C='Mod'; Height=0}
    @{Row=1022; B='lab.mod.index.title'; C='Detail modu'; Height=0}
    @{Row=1023; B='lab.mod.update.submit'; C='Aktualizovat'; Height=0}
    @{Row=1024; B='lab.mod.updated.message'; C='Mod [{{data.name}}] byl úspěšně aktualizován.'; Height=0}
)

foreach ($item in $data) {
    $r = $item.Row
    $srcRow = $r - 1
    $ws.Range("A" + $srcRow + ":C" + $srcRow).Copy() | Out-Null
    $ws.Range("A" + $r + ":C" + $r).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($r, 1).Value = "cs"
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    if ($item.Height -gt 0) {
        $ws.Rows.Item($r).RowHeight = $item.Height
    }
}

$excel.CutCopyMode = $false

$ws.Range("B1016").Select() | Out-Null
